# Update the ShipmentTracking numbers (column P) on rows 2-4 of the CRUD
# test-fixture worksheet. The new values are long, purely-numeric strings,
# so they must be forced to Text (otherwise Excel auto-converts them to a
# numeric value). We set the NumberFormat to Text ("@") before writing the
# value, then reset the cell style back to "Normal" so the text is stored
# without leaving a lingering custom number format / quote-prefix on the
# cell (matching how the cells were originally authored: plain string
# values with the default "General" style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "320018475104"
$ws.Range("P2").Style = "Normal"

$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Value = "320018514701"
$ws.Range("P3").Style = "Normal"

$ws.Range("P4").NumberFormat = "@"
$ws.Range("P4").Value = "320018475115"
$ws.Range("P4").Style = "Normal"
